$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Label" in H1, copying the formatting (style) from the
# neighboring header cell G1 (bold, bordered, centered)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# Fill in the new Label column (H) values for rows 2-11 and 12-21.
# Control patients (rows 2-6, 12-16) => 0 ; MDD patients (rows 7-11, 17-21) => 1
$labelValues = 0,0,0,0,0,1,1,1,1,1,0,0,0,0,0,1,1,1,1,1
for ($i = 0; $i -lt $labelValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labelValues[$i]
}
